$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from H1 (bold, bordered, centered style) onto the
# new header cells I1:J1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# New column data for I (I0) and J (IF), rows 2-22.
$data = @{
    2  = @(8, 8)
    3  = @(8, 9)
    4  = @(8, 8)
    5  = @(7, 7)
    6  = @(8, 8)
    7  = @(6, 6)
    8  = @(7, 7)
    9  = @(7, 8)
    10 = @(9, 9)
    11 = @(9, 9)
    12 = @(8, 8)
    13 = @(10, 10)
    14 = @(5, 5)
    15 = @(8, 8)
    16 = @(9, 9)
    17 = @(8, 8)
    18 = @(6, 6)
    19 = @(8, 8)
    20 = @(7, 7)
    21 = @(7, 7)
    22 = @(6, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value2 = $vals[0]
    $ws.Cells.Item($row, 10).Value2 = $vals[1]
}
